# Rename the worksheet from "Hoja1" to "prestaciones"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "prestaciones"
